$wb = $excel.ActiveWorkbook

# Rename the sheets (task order identifiers refreshed)
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687609378588"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168765139821"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687651413836"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687652014022"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687652879"

# Sheet1 - GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168760906989.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687609217362.csv"
$ws1.Range("B4").Value = "go_stims-1651168760922737.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687609368596.csv"

# Sheet2 - NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_6-16511687615548983.csv"
$ws2.Range("B3").Value = "OB-16511687628757625.csv"
$ws2.Range("B4").Value = "OB-16511687624972637.csv"
$ws2.Range("B5").Value = "TB-16511687650675337.csv"
$ws2.Range("B6").Value = "TB-16511687651216717.csv"
$ws2.Range("B7").Value = "ZB-match_5-16511687616140127.csv"
$ws2.Range("B8").Value = "OB-16511687622322848.csv"
$ws2.Range("B9").Value = "TB-16511687636134322.csv"
$ws2.Range("B10").Value = "ZB-match_7-16511687613688402.csv"

# Sheet4 - TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687651569788.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687651443772.csv"
$ws4.Range("B4").Value = "MM_stims-16511687651872.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687651579666.csv"
$ws4.Range("B6").Value = "MM_stims-16511687652003956.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687651872.csv"

# Sheet5 - vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511687652400582.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687652248745.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687652730813.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687652084358.csv"
